# Updates cryptocurrency price (column D) and 1h-volume-change (column E)
# values for the rows whose source data refreshed, per the commit diff.
# Values are plain text (matching the original t="inlineStr" cells), so a
# leading apostrophe forces Excel to keep numeric-looking strings (e.g.
# "290.22", "-8.98%") as text instead of auto-converting them to numbers;
# resetting the Style back to "Normal" afterwards drops the quote-prefix
# formatting flag that the apostrophe trick sets, so cell styling is
# untouched (matches the diff, which only changes text content).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'290.22"
$ws.Range("E2").Value = "'-8.98%"
$ws.Range("D2:E2").Style = "Normal"

$ws.Range("E3").Value = "'-2.47%"
$ws.Range("E3").Style = "Normal"

$ws.Range("D4").Value = "'5.053"
$ws.Range("E4").Value = "'-3.38%"
$ws.Range("D4:E4").Style = "Normal"

$ws.Range("D5").Value = "'0.07280"
$ws.Range("E5").Value = "'-5.32%"
$ws.Range("D5:E5").Style = "Normal"

$ws.Range("D6").Value = "'4.284"
$ws.Range("E6").Value = "'-1.50%"
$ws.Range("D6:E6").Style = "Normal"

$ws.Range("E7").Value = "'-6.98%"
$ws.Range("E7").Style = "Normal"

$ws.Range("D8").Value = "'0.9200"
$ws.Range("E8").Value = "'-1.92%"
$ws.Range("D8:E8").Style = "Normal"

$ws.Range("D9").Value = "'0.1150"
$ws.Range("E9").Value = "'-9.89%"
$ws.Range("D9:E9").Style = "Normal"

$ws.Range("D10").Value = "'0.1727"
$ws.Range("E10").Value = "'-6.32%"
$ws.Range("D10:E10").Style = "Normal"

$ws.Range("D11").Value = "'0.08644"
$ws.Range("E11").Value = "'-5.44%"
$ws.Range("D11:E11").Style = "Normal"

$ws.Range("D12").Value = "'0.04180"
$ws.Range("E12").Value = "'1.13%"
$ws.Range("D12:E12").Style = "Normal"

$ws.Range("E13").Value = "'0.17%"
$ws.Range("E13").Style = "Normal"

$ws.Range("D14").Value = "'0.001276"
$ws.Range("E14").Value = "'0.72%"
$ws.Range("D14:E14").Style = "Normal"

$ws.Range("D15").Value = "'0.005827"
$ws.Range("E15").Value = "'-1.23%"
$ws.Range("D15:E15").Style = "Normal"

$ws.Range("D16").Value = "'3.401"
$ws.Range("E16").Value = "'1.66%"
$ws.Range("D16:E16").Style = "Normal"

$ws.Range("D17").Value = "'2.336"
$ws.Range("E17").Value = "'-3.67%"
$ws.Range("D17:E17").Style = "Normal"

$ws.Range("D18").Value = "'0.3276"
$ws.Range("E18").Value = "'-2.24%"
$ws.Range("D18:E18").Style = "Normal"

$ws.Range("D19").Value = "'7.880"
$ws.Range("E19").Value = "'-6.24%"
$ws.Range("D19:E19").Style = "Normal"

$ws.Range("D20").Value = "'0.1342"
$ws.Range("E20").Value = "'-1.20%"
$ws.Range("D20:E20").Style = "Normal"

$ws.Range("D21").Value = "'0.2885"
$ws.Range("E21").Value = "'0.47%"
$ws.Range("D21:E21").Style = "Normal"

$ws.Range("D22").Value = "'0.03872"
$ws.Range("E22").Value = "'-4.25%"
$ws.Range("D22:E22").Style = "Normal"

$ws.Range("D23").Value = "'0.001270"
$ws.Range("E23").Value = "'-0.44%"
$ws.Range("D23:E23").Style = "Normal"

$ws.Range("D24").Value = "'0.003804"
$ws.Range("E24").Value = "'-7.36%"
$ws.Range("D24:E24").Style = "Normal"

$ws.Range("E25").Value = "'0.43%"
$ws.Range("E25").Style = "Normal"

$ws.Range("D26").Value = "'0.0003727"
$ws.Range("E26").Value = "'-95.02%"
$ws.Range("D26:E26").Style = "Normal"

$ws.Range("D38").Value = "'0.02315"
$ws.Range("E38").Value = "'-7.47%"
$ws.Range("D38:E38").Style = "Normal"

$ws.Range("D39").Value = "'0.04954"
$ws.Range("E39").Value = "'-6.33%"
$ws.Range("D39:E39").Style = "Normal"

$ws.Range("D40").Value = "'0.006690"
$ws.Range("E40").Value = "'224.61%"
$ws.Range("D40:E40").Style = "Normal"

$ws.Range("D41").Value = "'0.007679"
$ws.Range("E41").Value = "'-1.32%"
$ws.Range("D41:E41").Style = "Normal"

$ws.Range("D42").Value = "'0.1270"
$ws.Range("E42").Value = "'-2.80%"
$ws.Range("D42:E42").Style = "Normal"

$ws.Range("D43").Value = "'0.007391"
$ws.Range("E43").Value = "'4.54%"
$ws.Range("D43:E43").Style = "Normal"

$ws.Range("D44").Value = "'0.007074"
$ws.Range("E44").Value = "'-14.87%"
$ws.Range("D44:E44").Style = "Normal"

$ws.Range("D45").Value = "'0.2899"
$ws.Range("E45").Value = "'-16.42%"
$ws.Range("D45:E45").Style = "Normal"

$ws.Range("D46").Value = "'0.00006415"
$ws.Range("E46").Value = "'-4.07%"
$ws.Range("D46:E46").Style = "Normal"

$ws.Range("E47").Value = "'-0.39%"
$ws.Range("E47").Style = "Normal"

$ws.Range("D48").Value = "'0.02129"
$ws.Range("E48").Value = "'-89.12%"
$ws.Range("D48:E48").Style = "Normal"

$ws.Range("E49").Value = "'-0.50%"
$ws.Range("E49").Style = "Normal"

$ws.Range("E50").Value = "'-0.39%"
$ws.Range("E50").Style = "Normal"

$ws.Range("E51").Value = "'-0.39%"
$ws.Range("E51").Style = "Normal"
